# "Generate Report for Handback"
#
# The handback transform for the d5479062-... file failed because the
# handback package's file name (zvgmsof3.dwd) didn't match the handoff
# file name. Update the localization-status report to reflect this:
#   - flip the "Ready for handoff" status (row 3) to "Handback transform
#     failed" on the Overview sheet and on each per-language sheet
#   - record the mismatch message in the "Error Detail" column (P) of
#     the zh-cn and de-de sheets
#   - widen column P (Error Detail) on those sheets so the message is
#     readable

$wb = $excel.ActiveWorkbook

$statusText = "Handback transform failed"

$zhMessage = "Handback file name: zvgmsof3.dwd is different with handoff file name: d5479062-c2ec-43bd-b96c-87826f984d8e.f9ecb8a28f20e0a458cd8492fb3aa8b464733237.zh-cn."
$deMessage = "Handback file name: zvgmsof3.dwd is different with handoff file name: d5479062-c2ec-43bd-b96c-87826f984d8e.f9ecb8a28f20e0a458cd8492fb3aa8b464733237.de-de."

# --- Overview sheet: the zh-cn / de-de status cells for the d5479062 row ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = $statusText
$overview.Range("F3").Value = $statusText

# --- zh-cn sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $statusText
$zhcn.Range("P3").Value = $zhMessage
$zhcn.Columns.Item(16).ColumnWidth = 39.17

# --- de-de sheet ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $statusText
$dede.Range("P3").Value = $deMessage
$dede.Columns.Item(16).ColumnWidth = 39.17
